$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Dhh","Boc","ECs",3,1,3.659958333333333,10.979875,0.4781132044744068,0.4781132044744067,3,1,2.740639666666667,8.221919,0.05814502584416985,0.05814502584416984,10.03062698668056,90.27564288012499,0.02779990463060325,0.02779990463060324),
    @("ECs","Dhh","Boc","FAPs",3,1,3.659958333333333,10.979875,0.4781132044744068,0.4781132044744067,3,1,36.96800500000001,110.904015,0.7843079965148284,0.7843079965148283,135.3013579664584,1217.712221698125,0.3749880095086065,0.3749880095086063),
    @("ECs","Dhh","Boc","MuSCs",3,1,3.659958333333333,10.979875,0.4781132044744068,0.4781132044744067,3,1,7.396246333333333,22.188739,0.1569177223230415,0.1569177223230415,27.06995340306944,243.629580627625,0.07502443505869455,0.07502443505869452),
    @("ECs","Dhh","Boc","Resolving-Mac",3,1,3.659958333333333,10.979875,0.4781132044744068,0.4781132044744067,2,0.6666666666666666,0.02965966666666667,0.088979,0.0006292553179602461,0.0006292553179602461,0.1085531441805556,0.976978297625,0.0003008552765025351,0.0003008552765025349),
    @("FAPs","Dhh","Boc","ECs",3,1,2.444496,7.333488,0.3193330932870009,0.3193330932870008,3,1,2.740639666666667,8.221919,0.05814502584416985,0.05814502584416984,6.699482702608,60.295344323472,0.01856763096207137,0.01856763096207136),
    @("FAPs","Dhh","Boc","FAPs",3,1,2.444496,7.333488,0.3193330932870009,0.3193330932870008,3,1,36.96800500000001,110.904015,0.7843079965148284,0.7843079965148283,90.36814035048002,813.3132631543201,0.2504554986168104,0.2504554986168104),
    @("FAPs","Dhh","Boc","MuSCs",3,1,2.444496,7.333488,0.3193330932870009,0.3193330932870008,3,1,7.396246333333333,22.188739,0.1569177223230415,0.1569177223230415,18.080094576848,162.720851191632,0.05010902166096753,0.05010902166096751),
    @("FAPs","Dhh","Boc","Resolving-Mac",3,1,2.444496,7.333488,0.3193330932870009,0.3193330932870008,2,0.6666666666666666,0.02965966666666667,0.088979,0.0006292553179602461,0.0006292553179602461,0.07250293652800001,0.6525264287520001,0.0002009420471515407,0.0002009420471515406),
    @("MuSCs","Dhh","Boc","ECs",3,1,1.388093333333333,4.16428,0.1813315046964271,0.1813315046964271,3,1,2.740639666666667,8.221919,0.05814502584416985,0.05814502584416984,3.804263650368888,34.23837285331999,0.01054352502693596,0.01054352502693596),
    @("MuSCs","Dhh","Boc","FAPs",3,1,1.388093333333333,4.16428,0.1813315046964271,0.1813315046964271,3,1,36.96800500000001,110.904015,0.7843079965148284,0.7843079965148283,51.31504128713333,461.8353715842001,0.142219749153474,0.1422197491534739),
    @("MuSCs","Dhh","Boc","MuSCs",3,1,1.388093333333333,4.16428,0.1813315046964271,0.1813315046964271,3,1,7.396246333333333,22.188739,0.1569177223230415,0.1569177223230415,10.26668022699111,92.40012204291999,0.02845412670237325,0.02845412670237325),
    @("MuSCs","Dhh","Boc","Resolving-Mac",3,1,1.388093333333333,4.16428,0.1813315046964271,0.1813315046964271,2,0.6666666666666666,0.02965966666666667,0.088979,0.0006292553179602461,0.0006292553179602461,0.04117038556888888,0.37053347012,0.0001141038136439601,0.0001141038136439601),
    @("Resolving-Mac","Dhh","Boc","ECs",3,1,0.162456,0.487368,0.02122219754216535,0.02122219754216534,3,1,2.740639666666667,8.221919,0.05814502584416985,0.05814502584416984,0.445233357688,4.007100219192,0.001233965224559282,0.001233965224559282),
    @("Resolving-Mac","Dhh","Boc","FAPs",3,1,0.162456,0.487368,0.02122219754216535,0.02122219754216534,3,1,36.96800500000001,110.904015,0.7843079965148284,0.7843079965148283,6.005674220280001,54.05106798252001,0.01664473923593762,0.01664473923593761),
    @("Resolving-Mac","Dhh","Boc","MuSCs",3,1,0.162456,0.487368,0.02122219754216535,0.02122219754216534,3,1,7.396246333333333,22.188739,0.1569177223230415,0.1569177223230415,1.201564594328,10.814081348952,0.003330138901006236,0.003330138901006235),
    @("Resolving-Mac","Dhh","Boc","Resolving-Mac",3,1,0.162456,0.487368,0.02122219754216535,0.02122219754216534,2,0.6666666666666666,0.02965966666666667,0.088979,0.0006292553179602461,0.0006292553179602461,0.004818390808000001,0.04336551727200001,"1.335418066221041E-05","1.335418066221041E-05")
)


for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $v = $rowVals[$j]
        if ($j -ge 4 -and $v -is [string]) {
            $v = [double]$v
        }
        $ws.Cells.Item($rowNum, $colNum).Value2 = $v
    }
}
